$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the SearchURL cell (B2): append a trailing "?" and turn it into a
#    real hyperlink (this also applies the built-in "Hyperlink" style).
$searchUrl = "http://www.futureshop.ca/Search/SearchResults.aspx?"
$ws.Hyperlinks.Add($ws.Range("B2"), $searchUrl, "", "", $searchUrl)

# 2. Insert a brand-new "description / h4 / class / prod-title" block right
#    where the old "description / a / ..." block used to start (row 15),
#    pushing everything from there on down by one block (4 rows + 1 blank
#    separator = 5 rows).
$ws.Rows("15:19").Insert()
$ws.Range("A15").Value = "Item_Attribute_Name"
$ws.Range("B15").Value = "description"
$ws.Range("A16").Value = "html_tag"
$ws.Range("B16").Value = "h4"
$ws.Range("A17").Value = "html_tag_attribute_name"
$ws.Range("B17").Value = "class"
$ws.Range("A18").Value = "html_tag_attribute_value"
$ws.Range("B18").Value = "prod-title"

# 3. Insert two new "price" blocks (a "span / dollars" block, with its own
#    blank separator) right before the final "price / content_location"
#    block, which now lives 5 rows further down (old row 33 -> new row 38).
$ws.Rows("38:42").Insert()
$ws.Range("A38").Value = "Item_Attribute_Name"
$ws.Range("B38").Value = "price"
$ws.Range("A39").Value = "html_tag"
$ws.Range("B39").Value = "span"
$ws.Range("A40").Value = "html_tag_attribute_name"
$ws.Range("B40").Value = "class"
$ws.Range("A41").Value = "html_tag_attribute_value"
$ws.Range("B41").Value = "dollars"
